$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peru Liga 1")

# Row 325
$ws.Range("B325").Value = 7302796
$ws.Range("F325").Value = 'Sport Huancayo'
$ws.Range("G325").Value = 'Sport Boys'
$ws.Range("H325").Value = 1
$ws.Range("I325").Value = 0
$ws.Range("K325").Value = 1.727
$ws.Range("L325").Value = 3.75
$ws.Range("M325").Value = 4.333
$ws.Range("N325").Value = 1.25
$ws.Range("O325").Value = 5.25
$ws.Range("P325").Value = 10
$ws.Range("Q325").Value = -1.75
$ws.Range("R325").Value = 1.925
$ws.Range("S325").Value = 1.875
$ws.Range("T325").Value = 3
$ws.Range("U325").Value = 1.875
$ws.Range("V325").Value = 1.925
$ws.Range("W325").Value = 0.25
$ws.Range("Z325").Value = -1
$ws.Range("AA325").Value = 0.875
$ws.Range("AB325").Value = -1
$ws.Range("AC325").Value = 0.925

# Row 326
$ws.Range("B326").Value = 7302200
$ws.Range("F326").Value = 'Carlos Manucci'
$ws.Range("G326").Value = 'Deportivo Binacional'
$ws.Range("H326").Value = 3
$ws.Range("I326").Value = 2
$ws.Range("K326").Value = 2
$ws.Range("L326").Value = 3.2
$ws.Range("M326").Value = 3.75
$ws.Range("N326").Value = 1.75
$ws.Range("O326").Value = 3.4
$ws.Range("P326").Value = 4.333
$ws.Range("Q326").Value = -0.5
$ws.Range("R326").Value = 1.85
$ws.Range("S326").Value = 1.95
$ws.Range("T326").Value = 2.5
$ws.Range("U326").Value = 1.85
$ws.Range("V326").Value = 1.95
$ws.Range("W326").Value = 0.75
$ws.Range("Z326").Value = 0.8500000000000001
$ws.Range("AA326").Value = -1
$ws.Range("AB326").Value = 0.8500000000000001
$ws.Range("AC326").Value = -1

# Row 332
$ws.Range("B332").Value = 7384628
$ws.Range("F332").Value = 'Deportivo Binacional'
$ws.Range("G332").Value = 'FBC Melgar'
$ws.Range("H332").Value = 1
$ws.Range("I332").Value = 2
$ws.Range("J332").Value = 'A'
$ws.Range("K332").Value = 2.75
$ws.Range("L332").Value = 3.3
$ws.Range("M332").Value = 2.375
$ws.Range("N332").Value = 3.3
$ws.Range("O332").Value = 3.6
$ws.Range("P332").Value = 2
$ws.Range("Q332").Value = 0.5
$ws.Range("U332").Value = 1.975
$ws.Range("V332").Value = 1.875
$ws.Range("W332").Value = -1
$ws.Range("Y332").Value = 1
$ws.Range("Z332").Value = -1
$ws.Range("AA332").Value = 1
$ws.Range("AB332").Value = 0.4875
$ws.Range("AC332").Value = -0.5

# Row 333
$ws.Range("B333").Value = 7384630
$ws.Range("F333").Value = 'Atletico Grau'
$ws.Range("G333").Value = 'Unin Comercio'
$ws.Range("H333").Value = 0
$ws.Range("I333").Value = 1
$ws.Range("K333").Value = 2.8
$ws.Range("L333").Value = 3.4
$ws.Range("M333").Value = 2.15
$ws.Range("N333").Value = 1.75
$ws.Range("P333").Value = 3.8
$ws.Range("Q333").Value = -0.75
$ws.Range("R333").Value = 2
$ws.Range("S333").Value = 1.8
$ws.Range("T333").Value = 3
$ws.Range("U333").Value = 1.85
$ws.Range("V333").Value = 1.95
$ws.Range("Y333").Value = 2.8
$ws.Range("AA333").Value = 0.8
$ws.Range("AB333").Value = -1
$ws.Range("AC333").Value = 0.95

# Row 334
$ws.Range("B334").Value = 7384627
$ws.Range("F334").Value = 'Universitario de Deportes'
$ws.Range("G334").Value = 'Sport Huancayo'
$ws.Range("H334").Value = 2
$ws.Range("I334").Value = 0
$ws.Range("J334").Value = 'H'
$ws.Range("K334").Value = 1.25
$ws.Range("L334").Value = 5
$ws.Range("M334").Value = 12
$ws.Range("N334").Value = 1.181
$ws.Range("O334").Value = 6
$ws.Range("P334").Value = 13
$ws.Range("Q334").Value = -1.75
$ws.Range("R334").Value = 1.8
$ws.Range("S334").Value = 2
$ws.Range("T334").Value = 2.75
$ws.Range("W334").Value = 0.181
$ws.Range("Y334").Value = -1
$ws.Range("Z334").Value = 0.4
$ws.Range("AA334").Value = -0.5

# Row 335
$ws.Range("B335").Value = 7384626
$ws.Range("F335").Value = 'Sporting Cristal'
$ws.Range("G335").Value = 'Alianza Atletico'
$ws.Range("H335").Value = 3
$ws.Range("I335").Value = 0
$ws.Range("J335").Value = 'H'
$ws.Range("K335").Value = 1.3
$ws.Range("L335").Value = 5
$ws.Range("M335").Value = 9
$ws.Range("N335").Value = 1.166
$ws.Range("O335").Value = 6.5
$ws.Range("P335").Value = 13
$ws.Range("Q335").Value = -2
$ws.Range("R335").Value = 1.85
$ws.Range("S335").Value = 1.95
$ws.Range("T335").Value = 3.25
$ws.Range("U335").Value = 2
$ws.Range("V335").Value = 1.8
$ws.Range("W335").Value = 0.1659999999999999
$ws.Range("Y335").Value = -1
$ws.Range("Z335").Value = 0.8500000000000001
$ws.Range("AA335").Value = -1
$ws.Range("AB335").Value = -0.5
$ws.Range("AC335").Value = 0.4

# Row 337
$ws.Range("B337").Value = 7384629
$ws.Range("F337").Value = 'Deportivo Garcilaso'
$ws.Range("G337").Value = 'Alianza Lima'
$ws.Range("H337").Value = 0
$ws.Range("I337").Value = 1
$ws.Range("J337").Value = 'A'
$ws.Range("K337").Value = 2.625
$ws.Range("L337").Value = 3.3
$ws.Range("M337").Value = 2.5
$ws.Range("N337").Value = 2.7
$ws.Range("O337").Value = 3.4
$ws.Range("P337").Value = 2.375
$ws.Range("Q337").Value = 0
$ws.Range("R337").Value = 2.025
$ws.Range("S337").Value = 1.775
$ws.Range("T337").Value = 2.25
$ws.Range("U337").Value = 1.825
$ws.Range("V337").Value = 1.975
$ws.Range("W337").Value = -1
$ws.Range("Y337").Value = 1.375
$ws.Range("Z337").Value = -1
$ws.Range("AA337").Value = 0.7749999999999999
$ws.Range("AB337").Value = -1
$ws.Range("AC337").Value = 0.9750000000000001

# Row 377
$ws.Range("B377").Value = 7818817
$ws.Range("F377").Value = 'Sport Boys'
$ws.Range("G377").Value = 'Cusco FC'
$ws.Range("K377").Value = 2.2
$ws.Range("L377").Value = 3.2
$ws.Range("M377").Value = 3.2
$ws.Range("N377").Value = 1.7
$ws.Range("O377").Value = 3.6
$ws.Range("P377").Value = 5
$ws.Range("Q377").Value = -0.75
$ws.Range("R377").Value = 1.9
$ws.Range("S377").Value = 1.95
$ws.Range("T377").Value = 2.5
$ws.Range("U377").Value = 1.975
$ws.Range("V377").Value = 1.875

# Row 378
$ws.Range("B378").Value = 7818816
$ws.Range("F378").Value = 'UTC Cajamarca'
$ws.Range("G378").Value = 'Universitario de Deportes'
$ws.Range("K378").Value = 3.3
$ws.Range("L378").Value = 3.3
$ws.Range("M378").Value = 2.1
$ws.Range("N378").Value = 4.333
$ws.Range("O378").Value = 3.1
$ws.Range("P378").Value = 1.95
$ws.Range("Q378").Value = 0.5
$ws.Range("R378").Value = 1.85
$ws.Range("S378").Value = 2
$ws.Range("T378").Value = 2
$ws.Range("U378").Value = 1.8
$ws.Range("V378").Value = 2.05

# Row 379
$ws.Range("N379").Value = 1.166
$ws.Range("O379").Value = 6.5
$ws.Range("P379").Value = 19
$ws.Range("R379").Value = 1.825
$ws.Range("S379").Value = 2.025
$ws.Range("T379").Value = 3.25
$ws.Range("U379").Value = 2.025
$ws.Range("V379").Value = 1.825

# Row 380
$ws.Range("N380").Value = 4.75
$ws.Range("O380").Value = 3.8
$ws.Range("P380").Value = 1.666
$ws.Range("Q380").Value = 0.75
$ws.Range("R380").Value = 1.975
$ws.Range("S380").Value = 1.875
$ws.Range("U380").Value = 2
$ws.Range("V380").Value = 1.85

# Row 382
$ws.Range("N382").Value = 1.2
$ws.Range("O382").Value = 6.5
$ws.Range("P382").Value = 15
$ws.Range("Q382").Value = -2
$ws.Range("T382").Value = 3.25
$ws.Range("U382").Value = 2.05
$ws.Range("V382").Value = 1.8
